# Fruta / hortaliza, semanal
#
# This workbook holds a weekly price list for "Durazno" (peach) at
# "Terminal La Palmera de La Serena". A new week's worth of rows for the
# variety "Kakamas" is inserted at the top of the data block (the data
# block begins at row 368), pushing every existing data row down by 3
# rows and growing the used range from A1:T433 to A1:T436.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the existing row 368 (Excel shifts
# rows 368:433 down to 371:436, carrying their formatting/styles with
# them - this matches the diff exactly, where every old row N (368..433)
# reappears, unchanged, as new row N+3).
$ws.Rows("368:370").Insert()

# Fill in the 3 newly-inserted rows (368, 369, 370) with the new
# "Kakamas" price entries. Columns A,B,C,E,F,G,H,I,J,Q,R are identical to
# the rest of the block and are copied from the shifted row 371 (the old
# row 368, "Carson / Especial").
$commonA = 8
$commonB = "Terminal La Palmera de La Serena"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103004
$commonJ = "Durazno"
$commonQ = "`$/bins (400 kilos)"
$commonR = "Región de O'Higgins"
$commonT = 400
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(44637)

$rows = @(
    @{ Row = 368; K = "Kakamas"; L = "Especial"; M = 20; N = 445000; O = 450000; P = 447500; S = 1119 },
    @{ Row = 369; K = "Kakamas"; L = "Primera";  M = 20; N = 405000; O = 410000; P = 407500; S = 1019 },
    @{ Row = 370; K = "Kakamas"; L = "Segunda";  M = 16; N = 355000; O = 360000; P = 357500; S = 894  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $newDate
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $commonQ
    $ws.Cells.Item($row, 18).Value = $commonR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $commonT
}
